# Automatic update of files.
# Reassign the per-record fields (Id, Ost/East, Nord/North, Publik kommentar)
# across rows 3-23, and align the empty Alder-Stadium/Kon/Aktivitet/Metod
# marker cells (K:N) with the relocated record in rows 4 and 8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 108718934
$ws.Range("Q3").Value = 437251.0097213153
$ws.Range("R3").Value = 7006251.485929966
$ws.Range("AC3").Value = ""

# Row 4
$ws.Range("A4").Value = 108718930
$ws.Range("Q4").Value = 437403.8165098006
$ws.Range("R4").Value = 7006294.143844516
$ws.Range("AC4").Value = 'ringhack'
$ws.Range("K4").NumberFormat = "General"
$ws.Range("L4").NumberFormat = "General"
$ws.Range("M4").NumberFormat = "General"
$ws.Range("N4").NumberFormat = "General"

# Row 5
$ws.Range("A5").Value = 108718928
$ws.Range("Q5").Value = 437260.4723100049
$ws.Range("R5").Value = 7006319.508680805
$ws.Range("AC5").Value = 'ringhack äldre'

# Row 6
$ws.Range("A6").Value = 108718932
$ws.Range("Q6").Value = 437306.259954496
$ws.Range("R6").Value = 7006254.930521684
$ws.Range("AC6").Value = 'ringhack äldre'

# Row 7
$ws.Range("A7").Value = 108718925
$ws.Range("Q7").Value = 437441.8026971049
$ws.Range("R7").Value = 7006433.882732502
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 108718923
$ws.Range("Q8").Value = 437693.1568600214
$ws.Range("R8").Value = 7006333.720904024
$ws.Range("AC8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""

# Row 9
$ws.Range("A9").Value = 108718940
$ws.Range("Q9").Value = 437142.4856154905
$ws.Range("R9").Value = 7006787.91673798
$ws.Range("AC9").Value = ""

# Row 10
$ws.Range("A10").Value = 108718927
$ws.Range("Q10").Value = 437347.3375837047
$ws.Range("R10").Value = 7006460.555242279
$ws.Range("AC10").Value = ""

# Row 11
$ws.Range("A11").Value = 108718949
$ws.Range("Q11").Value = 437232.1395232935
$ws.Range("R11").Value = 7006653.843598615
$ws.Range("AC11").Value = 'ringhack'

# Row 12
$ws.Range("A12").Value = 108718943
$ws.Range("Q12").Value = 437156.0212149328
$ws.Range("R12").Value = 7006739.780347738
$ws.Range("AC12").Value = 'ringhack'

# Row 13
$ws.Range("A13").Value = 108718947
$ws.Range("Q13").Value = 437196.0548698546
$ws.Range("R13").Value = 7006728.614212831
$ws.Range("AC13").Value = 'ringhack äldre'

# Row 14
$ws.Range("A14").Value = 108718926
$ws.Range("Q14").Value = 437430.2372204551
$ws.Range("R14").Value = 7006420.556342849
$ws.Range("AC14").Value = 'ringhack'

# Row 15
$ws.Range("A15").Value = 108718953
$ws.Range("Q15").Value = 438033.4045416421
$ws.Range("R15").Value = 7006292.397694888
$ws.Range("AC15").Value = 'ringhack färska'

# Row 16
$ws.Range("A16").Value = 108718948
$ws.Range("Q16").Value = 437215.8330926147
$ws.Range("R16").Value = 7006675.839449953
$ws.Range("AC16").Value = 'ringhack äldre'

# Row 17
$ws.Range("A17").Value = 108718938
$ws.Range("Q17").Value = 437083.5155771806
$ws.Range("R17").Value = 7006709.125224494
$ws.Range("AC17").Value = 'ringhack'

# Row 18
$ws.Range("A18").Value = 108718939
$ws.Range("Q18").Value = 437092.8160429197
$ws.Range("R18").Value = 7006699.008007247
$ws.Range("AC18").Value = 'ringhack färska'

# Row 19
$ws.Range("A19").Value = 108718950
$ws.Range("Q19").Value = 437330.4075466889
$ws.Range("R19").Value = 7006637.031554679
$ws.Range("AC19").Value = 'ringhack äldre'

# Row 20
$ws.Range("A20").Value = 108718941
$ws.Range("Q20").Value = 437149.3458266784
$ws.Range("R20").Value = 7006791.847929343
$ws.Range("AC20").Value = 'ringhack'

# Row 21
$ws.Range("A21").Value = 108718944
$ws.Range("Q21").Value = 437153.3087088031
$ws.Range("R21").Value = 7006739.833113052
$ws.Range("AC21").Value = 'ringhack färska'

# Row 22
$ws.Range("A22").Value = 108718952
$ws.Range("Q22").Value = 437864.8221849522
$ws.Range("R22").Value = 7006346.224326964
$ws.Range("AC22").Value = ""

# Row 23
$ws.Range("A23").Value = 108718946
$ws.Range("Q23").Value = 437193.3774615529
$ws.Range("R23").Value = 7006730.472810662
$ws.Range("AC23").Value = 'ringhack'

